$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3057.7334
$ws.Range("J40").Value = 6499.5
$ws.Range("L40").Value = 6499.5
$ws.Range("N40").Value = -6849.5
$ws.Range("H43").Value = 14500
$ws.Range("J43").Value = 14500
$ws.Range("L43").Value = 14500
$ws.Range("N43").Value = -14638
$ws.Range("H86").Value = 50977.555
$ws.Range("I86").Value = 4325
$ws.Range("J86").Value = 88299.60000000001
$ws.Range("K86").Value = 4325
$ws.Range("L86").Value = 88299.60000000001
$ws.Range("M86").Value = -3202
$ws.Range("N86").Value = -90545.60000000001
$ws.Range("H89").Value = 50977.555
$ws.Range("I89").Value = 4325
$ws.Range("J89").Value = 88299.60000000001
$ws.Range("K89").Value = 21625
$ws.Range("L89").Value = 441498
$ws.Range("M89").Value = -16009
$ws.Range("N89").Value = -452730
$ws.Range("H107").Value = 6516.609
$ws.Range("I107").Value = 8869.200000000001
$ws.Range("K107").Value = 8869.200000000001
$ws.Range("M107").Value = -6949.200000000001
$ws.Range("H132").Value = 5725.1943
$ws.Range("I132").Value = 5698.4194
$ws.Range("K132").Value = 17095.2582
$ws.Range("M132").Value = -14565.2582
$ws.Range("H138").Value = 288865.16
$ws.Range("I138").Value = 527500.2
$ws.Range("J138").Value = 4775.857
$ws.Range("K138").Value = 1582500.6
$ws.Range("L138").Value = 14327.571
$ws.Range("M138").Value = -1577360.6
$ws.Range("N138").Value = -24607.571

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 2004.1052
$ws.Range("I110").Value = 1258.1538
$ws.Range("K110").Value = 1258.1538
$ws.Range("M110").Value = 786.8462
$ws.Range("H122").Value = 393817.97
$ws.Range("I122").Value = 4050
$ws.Range("K122").Value = 12150
$ws.Range("M122").Value = -9700
$ws.Range("H132").Value = 3393.3333
$ws.Range("I132").Value = 2448.5
$ws.Range("J132").Value = 5719.077
$ws.Range("K132").Value = 7345.5
$ws.Range("L132").Value = 17157.231
$ws.Range("M132").Value = -4815.5
$ws.Range("N132").Value = -22217.231

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1454.4546
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 1299.9
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 1299.9
$ws.Range("M80").Value = -2002
$ws.Range("N80").Value = -3295.9
$ws.Range("H83").Value = 1454.4546
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 1299.9
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 6499.5
$ws.Range("M83").Value = -10008
$ws.Range("N83").Value = -16483.5
$ws.Range("H86").Value = 4801.154
$ws.Range("I86").Value = 7318.5713
$ws.Range("J86").Value = 1864.1666
$ws.Range("K86").Value = 7318.5713
$ws.Range("L86").Value = 1864.1666
$ws.Range("M86").Value = -6195.5713
$ws.Range("N86").Value = -4110.1666
$ws.Range("H89").Value = 4801.154
$ws.Range("I89").Value = 7318.5713
$ws.Range("J89").Value = 1864.1666
$ws.Range("K89").Value = 36592.85649999999
$ws.Range("L89").Value = 9320.833000000001
$ws.Range("M89").Value = -30976.85649999999
$ws.Range("N89").Value = -20552.833
$ws.Range("H94").Value = 9498.031999999999
$ws.Range("I94").Value = 12068.542
$ws.Range("J94").Value = 684.8570999999999
$ws.Range("K94").Value = 12068.542
$ws.Range("L94").Value = 684.8570999999999
$ws.Range("M94").Value = -11617.542
$ws.Range("N94").Value = -1586.8571
$ws.Range("H105").Value = 105489
$ws.Range("I105").Value = 252225
$ws.Range("K105").Value = 252225
$ws.Range("M105").Value = -250478
$ws.Range("H107").Value = 1992.6
$ws.Range("I107").Value = 2353.5454
$ws.Range("K107").Value = 2353.5454
$ws.Range("M107").Value = -433.5454

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 853.1667
$ws.Range("I16").Value = 853.1667
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 853.1667
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -566.1667
$ws.Range("N16").ClearContents()
$ws.Range("H99").Value = 6836380.5
$ws.Range("I99").Value = 16591960
$ws.Range("K99").Value = 16591960
$ws.Range("M99").Value = -16590462
$ws.Range("H105").Value = 424800
$ws.Range("I105").Value = 529750
$ws.Range("J105").Value = 5000
$ws.Range("K105").Value = 529750
$ws.Range("L105").Value = 5000
$ws.Range("M105").Value = -528003
$ws.Range("N105").Value = -8494
$ws.Range("H113").Value = 853.1667
$ws.Range("I113").Value = 853.1667
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 853.1667
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1316.8333
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 11418.75
$ws.Range("I122").Value = 30612.5
$ws.Range("K122").Value = 91837.5
$ws.Range("M122").Value = -89387.5
$ws.Range("H126").Value = 6836380.5
$ws.Range("I126").Value = 16591960
$ws.Range("K126").Value = 49775880
$ws.Range("M126").Value = -49773410

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 4488.222
$ws.Range("I137").Value = 1361.9166
$ws.Range("K137").Value = 4085.7498
$ws.Range("M137").Value = 1014.2502

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 173.5
$ws.Range("I2").Value = 148.23529
$ws.Range("K2").Value = 148.23529
$ws.Range("M2").Value = -35.23528999999999
$ws.Range("H107").Value = 378.8
$ws.Range("I107").Value = 421.5
$ws.Range("J107").Value = 314.75
$ws.Range("K107").Value = 421.5
$ws.Range("L107").Value = 314.75
$ws.Range("M107").Value = 1498.5
$ws.Range("N107").Value = -4154.75

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 8738.308000000001
$ws.Range("I22").Value = 15931.077
$ws.Range("J22").Value = 1545.5385
$ws.Range("K22").Value = 15931.077
$ws.Range("L22").Value = 1545.5385
$ws.Range("M22").Value = -15636.077
$ws.Range("N22").Value = -2135.5385
$ws.Range("H27").Value = 8738.308000000001
$ws.Range("I27").Value = 15931.077
$ws.Range("J27").Value = 1545.5385
$ws.Range("K27").Value = 15931.077
$ws.Range("L27").Value = 1545.5385
$ws.Range("M27").Value = -15824.077
$ws.Range("N27").Value = -1759.5385
$ws.Range("H46").Value = 4767.909
$ws.Range("J46").Value = 5666.3335
$ws.Range("L46").Value = 5666.3335
$ws.Range("N46").Value = -6042.3335
$ws.Range("H53").Value = 12000
$ws.Range("J53").Value = 12000
$ws.Range("L53").Value = 12000
$ws.Range("N53").Value = -13036
$ws.Range("H55").Value = 1918.9
$ws.Range("I55").Value = 625
$ws.Range("J55").Value = 2242.375
$ws.Range("K55").Value = 625
$ws.Range("L55").Value = 2242.375
$ws.Range("M55").Value = -452
$ws.Range("N55").Value = -2588.375
$ws.Range("H82").Value = 3500.3845
$ws.Range("J82").Value = 3125
$ws.Range("L82").Value = 3125
$ws.Range("N82").Value = -3847
$ws.Range("H85").Value = 3500.3845
$ws.Range("J85").Value = 3125
$ws.Range("L85").Value = 3125
$ws.Range("N85").Value = -5621
$ws.Range("H136").Value = 3437.1365
$ws.Range("I136").Value = 2727.2104
$ws.Range("K136").Value = 8181.6312
$ws.Range("M136").Value = -5631.6312

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 12889
$ws.Range("I2").Value = 7240.5454
$ws.Range("J2").Value = 33600
$ws.Range("K2").Value = 7240.5454
$ws.Range("L2").Value = 33600
$ws.Range("M2").Value = -7128.5454
$ws.Range("N2").Value = -33824
$ws.Range("H4").Value = 676.9231
$ws.Range("I4").Value = 333.33334
$ws.Range("J4").Value = 1450
$ws.Range("K4").Value = 333.33334
$ws.Range("L4").Value = 1450
$ws.Range("M4").Value = -1676

Write-Host "Applied all cell updates."